$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value for every data row (2-218).
# All of these change from 45188 (2023-09-19) to 45189 (2023-09-20).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 218 }

$ws.Range("C2:C$lastRow").Value = 45189
